$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data table (header + 13 data rows) to reflect new mock data set
$ws.Cells.Item(1, 1).Value = "CLINIC"
$ws.Cells.Item(1, 2).Value = "RESPONSE"
$ws.Cells.Item(1, 3).Value = "COMMENTS"

$ws.Cells.Item(2, 1).Value = "Theatre Treatment Suite Implants"
$ws.Cells.Item(2, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(2, 3).Value = "Long wait times"

$ws.Cells.Item(3, 1).Value = "Labour and Delivery Suite"
$ws.Cells.Item(3, 2).Value = "Unlikely"
$ws.Cells.Item(3, 3).Value = "Long wait times"

$ws.Cells.Item(4, 1).Value = "Sitwell"
$ws.Cells.Item(4, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(4, 3).Value = "Food was terrible"

$ws.Cells.Item(5, 1).Value = "A&E"
$ws.Cells.Item(5, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(5, 3).Value = "Waited for long time for poor service"

$ws.Cells.Item(6, 1).Value = "Gynaecology"
$ws.Cells.Item(6, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(6, 3).Value = "Waited for long time for poor service"

$ws.Cells.Item(7, 1).Value = "Bone Health"
$ws.Cells.Item(7, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(7, 3).Value = "Waited for long time for poor service"

$ws.Cells.Item(8, 1).Value = "Bone Health"
$ws.Cells.Item(8, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(8, 3).Value = "staff was rude"

$ws.Cells.Item(9, 1).Value = "Radiology"
$ws.Cells.Item(9, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(9, 3).Value = "Doctors are patronising and made me feel bad"

$ws.Cells.Item(10, 1).Value = "Labour and Delivery Suite"
$ws.Cells.Item(10, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(10, 3).Value = "staff tried to deal with me quickly rather than correctly. Not appropriate and i shouldve have been taken care of better. Would not recommend."

$ws.Cells.Item(11, 1).Value = "Sitwell"
$ws.Cells.Item(11, 2).Value = "Unlikely"
$ws.Cells.Item(11, 3).Value = "staff tried to deal with me quickly rather than correctly. Not appropriate and i shouldve have been taken care of better. Would not recommend."

$ws.Cells.Item(12, 1).Value = "Rehab Services"
$ws.Cells.Item(12, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(12, 3).Value = "Waited over 5 hours"

$ws.Cells.Item(13, 1).Value = "Radiology"
$ws.Cells.Item(13, 2).Value = "Extremely Unlikely"
$ws.Cells.Item(13, 3).Value = "Waited too long to find a parking spot"

$ws.Cells.Item(14, 1).Value = "Theatre Treatment Suite Implants"
$ws.Cells.Item(14, 2).Value = "Unlikely"
$ws.Cells.Item(14, 3).Value = "Felt as if i was not a priority"

# Update selection to match the saved view state (A6:XFD6, i.e. row 6 selected)
$ws.Range("A6:XFD6").Select()
